$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.02
$ws.Range("L2").Value = 1.25
$ws.Range("N2").Value = 6.4
$ws.Range("P2").Value = 2.88
$ws.Range("W2").Value = 1.98
$ws.Range("F3").Value = 5.2
$ws.Range("G3").Value = 6.2
$ws.Range("H3").Value = 1.56
$ws.Range("N3").Value = 5.9
$ws.Range("F5").Value = 3.2
$ws.Range("H5").Value = 2.2
$ws.Range("K5").Value = 3.9
$ws.Range("P5").Value = 1.92
$ws.Range("Q5").Value = 1.87
$ws.Range("P6").Value = 2.06
$ws.Range("Q6").Value = 1.74
$ws.Range("F7").Value = 1.96
$ws.Range("G7").Value = 2.2
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3.9
$ws.Range("K7").Value = 4.7
$ws.Range("P7").Value = 2.52
$ws.Range("F8").Value = 2.26
$ws.Range("K8").Value = 4.2
$ws.Range("F9").Value = 6.2
$ws.Range("G9").Value = 8.800000000000001
$ws.Range("H9").Value = 1.48
$ws.Range("I9").Value = 1.58
$ws.Range("J9").Value = 4.5
$ws.Range("P9").Value = 2.4
$ws.Range("Q9").Value = 1.57
$ws.Range("G10").Value = 2.38
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3.8
$ws.Range("J10").Value = 3.55
$ws.Range("P10").Value = 1.95
$ws.Range("F13").Value = 5.2
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 1.63
$ws.Range("I13").Value = 1.78
$ws.Range("J13").Value = 3.85
$ws.Range("F14").Value = 2.52
$ws.Range("G14").Value = 2.9
$ws.Range("I14").Value = 3.4
$ws.Range("J14").Value = 3.15
$ws.Range("G15").Value = 6.2
$ws.Range("H15").Value = 1.56
$ws.Range("I15").Value = 1.7
$ws.Range("F16").Value = 3.3
$ws.Range("I16").Value = 2
$ws.Range("P16").Value = 2.06
$ws.Range("Q16").Value = 1.57
$ws.Range("H17").Value = 1.68
$ws.Range("G18").Value = 1.29
$ws.Range("I18").Value = 14.5
$ws.Range("J18").Value = 6.8
$ws.Range("K18").Value = 8.4
$ws.Range("N18").Value = 8.4
$ws.Range("O18").Value = 1.1
$ws.Range("P18").Value = 3.45
$ws.Range("Q18").Value = 1.32
$ws.Range("R18").Value = 2
$ws.Range("S18").Value = 1.8
$ws.Range("T18").Value = 1.68
$ws.Range("U18").Value = 2.18
$ws.Range("X18").Value = 55
$ws.Range("Z18").Value = 180
$ws.Range("AA18").Value = 530
$ws.Range("AB18").Value = 16
$ws.Range("AC18").Value = 19
$ws.Range("AD18").Value = 46
$ws.Range("AE18").Value = 180
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 29
$ws.Range("AI18").Value = 130
$ws.Range("AJ18").Value = 13
$ws.Range("AK18").Value = 15
$ws.Range("AL18").Value = 28
$ws.Range("AN18").Value = 3.35
$ws.Range("AO18").Value = 1000
$ws.Range("AB19").Value = 14
$ws.Range("AE19").Value = 29
$ws.Range("AI19").Value = 30
$ws.Range("AJ19").Value = 42
$ws.Range("AL19").Value = 28
$ws.Range("AN19").Value = 17.5
$ws.Range("AO19").Value = 19.5
